$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 55, shifting existing row 55 (and everything below) down by one.
$ws.Rows("55:55").Insert()

# Populate the new author row: Fabian Kellerer at IFIC (Instituto de Fisica Corpuscular).
$ws.Range("A55").Value = 'Kellerer'
$ws.Range("B55").Value = 'F.'
$ws.Range("E55").Value = 'Instituto de F\''isica Corpuscular (IFIC), CSIC \& Universitat de Val\`encia, Calle Catedr\''atico Jos\''e Beltr\''an, 2 '
$ws.Range("F55").Value = ' Paterna, E-46980, Spain'

# Match the author's final selection/scroll position on save.
$ws.Range("E55").Select()
